# Completes the "iden prob result" section on the weights sheet:
# adds +-cputime / +-gputime columns (E,F) next to the cpu_time/gpu_time
# table (rows 34-49) and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weights")

# --- header row for the per-sample table (row 34) ---
$ws.Range("E34").Value = "+-cputime"
$ws.Range("F34").Value = "+-gputime"

# --- per-sample absolute deviation formulas (rows 35-43) ---
$ws.Range("E35").Formula = '=ABS(C$47-C35)'
$ws.Range("F35").Formula = '=ABS(D$47-D35)'

$ws.Range("E36").Formula = '=ABS(C$47-C36)'
$ws.Range("F36").Formula = '=ABS(D$47-D36)'

$ws.Range("E37").Formula = '=ABS(C$47-C37)'
$ws.Range("F37").Formula = '=ABS(D$47-D37)'

$ws.Range("E38").Formula = '=ABS(C$48-C38)'
$ws.Range("F38").Formula = '=ABS(D$48-D38)'

$ws.Range("E39").Formula = '=ABS(C$48-C39)'
$ws.Range("F39").Formula = '=ABS(D$48-D39)'

$ws.Range("E40").Formula = '=ABS(C$48-C40)'
$ws.Range("F40").Formula = '=ABS(D$48-D40)'

$ws.Range("E41").Formula = '=ABS(C$49-C41)'
$ws.Range("F41").Formula = '=ABS(D$49-D41)'

$ws.Range("E42").Formula = '=ABS(C$49-C42)'
$ws.Range("F42").Formula = '=ABS(D$49-D42)'

$ws.Range("E43").Formula = '=ABS(C$49-C43)'
$ws.Range("F43").Formula = '=ABS(D$49-D43)'

# --- header row for the averages table (row 46) ---
$ws.Range("E46").Value = "+-cpu"
$ws.Range("F46").Value = "+-gpu"

# --- averages of the deviations (rows 47-49) ---
$ws.Range("E47").Formula = '=AVERAGE(E35:E37)'
$ws.Range("F47").Formula = '=AVERAGE(F35:F37)'

$ws.Range("E48").Formula = '=AVERAGE(E36:E38)'
$ws.Range("F48").Formula = '=AVERAGE(F36:F38)'

$ws.Range("E49").Formula = '=AVERAGE(E37:E39)'
$ws.Range("F49").Formula = '=AVERAGE(F37:F39)'

# --- move the view / active selection to the newly completed area ---
$ws.Activate()
$ws.Range("G47").Select()
